$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: add status comment
$ws.Range("H27").Value = "Set up the hybris server and done initialisation to set up the code.Development Done.Not Getting response from hybris sever"
$ws.Range("H27").WrapText = $true

# Row 14: add end date, status comment, and overall status
$ws.Range("F14").Value = 42222
$ws.Range("F14").NumberFormat = $ws.Range("F13").NumberFormat

$ws.Range("H14").Value = "Started Writing client code and methods.Waiting for updated webservice from Swarnima."
$ws.Range("H14").WrapText = $true

$ws.Range("I14").Value = "In progress"

# Row heights for the two updated rows
$ws.Rows.Item(14).RowHeight = 42
$ws.Rows.Item(27).RowHeight = 42

# Move the active selection to F18
$ws.Range("F18").Select()
